$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9106121063232422
$ws.Range("B1").Value = 1.724440097808838
$ws.Range("C1").Value = 4.722105979919434
$ws.Range("D1").Value = 4.885763645172119
$ws.Range("E1").Value = 1.685155987739563
